# Auto-generated script to apply 2023-09-24 crime data update
# Updates column J (year 2023 cumulative totals) across 47 worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 5590
$ws.Cells.Item(3, 10).Value = 5969
$ws.Cells.Item(4, 10).Value = 1295
$ws.Cells.Item(5, 10).Value = 458
$ws.Cells.Item(6, 10).Value = 7561
$ws.Cells.Item(7, 10).Value = 20873

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(3, 10).Value = 39
$ws.Cells.Item(6, 10).Value = 187
$ws.Cells.Item(7, 10).Value = 291

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 10).Value = 397
$ws.Cells.Item(6, 10).Value = 436

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 158
$ws.Cells.Item(6, 10).Value = 116
$ws.Cells.Item(7, 10).Value = 425

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 235
$ws.Cells.Item(3, 10).Value = 315
$ws.Cells.Item(6, 10).Value = 329
$ws.Cells.Item(7, 10).Value = 959

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 10).Value = 109
$ws.Cells.Item(7, 10).Value = 311

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 191
$ws.Cells.Item(7, 10).Value = 647

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(4, 10).Value = 24
$ws.Cells.Item(6, 10).Value = 186
$ws.Cells.Item(7, 10).Value = 530

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 91
$ws.Cells.Item(3, 10).Value = 129
$ws.Cells.Item(7, 10).Value = 325

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 10).Value = 90
$ws.Cells.Item(5, 10).Value = 65
$ws.Cells.Item(7, 10).Value = 612
$ws.Cells.Item(10, 10).Value = 142
$ws.Cells.Item(15, 10).Value = 228
$ws.Cells.Item(19, 10).Value = 609
$ws.Cells.Item(20, 10).Value = 432
$ws.Cells.Item(21, 10).Value = 59
$ws.Cells.Item(27, 10).Value = 126
$ws.Cells.Item(29, 10).Value = 1173
$ws.Cells.Item(31, 10).Value = 190
$ws.Cells.Item(33, 10).Value = 959
$ws.Cells.Item(37, 10).Value = 647
$ws.Cells.Item(42, 10).Value = 867
$ws.Cells.Item(45, 10).Value = 30
$ws.Cells.Item(46, 10).Value = 70
$ws.Cells.Item(47, 10).Value = 159
$ws.Cells.Item(48, 10).Value = 243
$ws.Cells.Item(52, 10).Value = 526
$ws.Cells.Item(53, 10).Value = 291
$ws.Cells.Item(54, 10).Value = 406
$ws.Cells.Item(57, 10).Value = 87
$ws.Cells.Item(60, 10).Value = 128
$ws.Cells.Item(63, 10).Value = 79
$ws.Cells.Item(65, 10).Value = 530
$ws.Cells.Item(67, 10).Value = 794
$ws.Cells.Item(72, 10).Value = 86
$ws.Cells.Item(76, 10).Value = 304
$ws.Cells.Item(78, 10).Value = 257
$ws.Cells.Item(79, 10).Value = 596
$ws.Cells.Item(80, 10).Value = 32
$ws.Cells.Item(82, 10).Value = 27
$ws.Cells.Item(83, 10).Value = 425
$ws.Cells.Item(84, 10).Value = 179
$ws.Cells.Item(85, 10).Value = 876
$ws.Cells.Item(86, 10).Value = 127
$ws.Cells.Item(87, 10).Value = 72
$ws.Cells.Item(88, 10).Value = 224
$ws.Cells.Item(89, 10).Value = 274
$ws.Cells.Item(91, 10).Value = 232
$ws.Cells.Item(94, 10).Value = 209
$ws.Cells.Item(95, 10).Value = 311
$ws.Cells.Item(97, 10).Value = 170
$ws.Cells.Item(98, 10).Value = 151
$ws.Cells.Item(99, 10).Value = 325
$ws.Cells.Item(101, 10).Value = 20873

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 10).Value = 73
$ws.Cells.Item(6, 10).Value = 52
$ws.Cells.Item(7, 10).Value = 190

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 10).Value = 197
$ws.Cells.Item(3, 10).Value = 303
$ws.Cells.Item(4, 10).Value = 61
$ws.Cells.Item(6, 10).Value = 211
$ws.Cells.Item(7, 10).Value = 794

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 10).Value = 54
$ws.Cells.Item(7, 10).Value = 179

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 10).Value = 80
$ws.Cells.Item(7, 10).Value = 406

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 352
$ws.Cells.Item(3, 10).Value = 406
$ws.Cells.Item(6, 10).Value = 306
$ws.Cells.Item(7, 10).Value = 1173

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 10).Value = 38
$ws.Cells.Item(6, 10).Value = 123
$ws.Cells.Item(7, 10).Value = 243

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 10).Value = 152
$ws.Cells.Item(3, 10).Value = 178
$ws.Cells.Item(6, 10).Value = 226
$ws.Cells.Item(7, 10).Value = 609

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 10).Value = 63
$ws.Cells.Item(6, 10).Value = 169
$ws.Cells.Item(7, 10).Value = 304

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 190
$ws.Cells.Item(3, 10).Value = 174
$ws.Cells.Item(6, 10).Value = 446
$ws.Cells.Item(7, 10).Value = 867

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 10).Value = 31
$ws.Cells.Item(6, 10).Value = 78
$ws.Cells.Item(7, 10).Value = 142

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 10).Value = 72
$ws.Cells.Item(7, 10).Value = 257

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(2, 10).Value = 21
$ws.Cells.Item(6, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 70

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 10).Value = 54
$ws.Cells.Item(7, 10).Value = 232

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 59

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 10).Value = 209
$ws.Cells.Item(6, 10).Value = 169
$ws.Cells.Item(7, 10).Value = 596

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 10).Value = 151
$ws.Cells.Item(6, 10).Value = 113
$ws.Cells.Item(7, 10).Value = 432

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 188
$ws.Cells.Item(3, 10).Value = 184
$ws.Cells.Item(5, 10).Value = 17
$ws.Cells.Item(7, 10).Value = 612

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 10).Value = 37
$ws.Cells.Item(6, 10).Value = 114
$ws.Cells.Item(7, 10).Value = 209

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(6, 10).Value = 75
$ws.Cells.Item(7, 10).Value = 159

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 10).Value = 66
$ws.Cells.Item(6, 10).Value = 95
$ws.Cells.Item(7, 10).Value = 228

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(2, 10).Value = 26
$ws.Cells.Item(6, 10).Value = 94
$ws.Cells.Item(7, 10).Value = 151

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 10).Value = 117
$ws.Cells.Item(7, 10).Value = 170

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 224

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 10).Value = 86
$ws.Cells.Item(3, 10).Value = 75
$ws.Cells.Item(6, 10).Value = 82
$ws.Cells.Item(7, 10).Value = 274

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(7, 10).Value = 65

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(6, 10).Value = 44
$ws.Cells.Item(7, 10).Value = 126

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(6, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 127

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 87

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(2, 10).Value = 46
$ws.Cells.Item(7, 10).Value = 128

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 229
$ws.Cells.Item(6, 10).Value = 255
$ws.Cells.Item(7, 10).Value = 876

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(3, 10).Value = 26
$ws.Cells.Item(4, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 86

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(5, 10).Value = 18
$ws.Cells.Item(6, 10).Value = 27

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Cells.Item(6, 10).Value = 10
$ws.Cells.Item(7, 10).Value = 30

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 10).Value = 162
$ws.Cells.Item(7, 10).Value = 526

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(6, 10).Value = 34
$ws.Cells.Item(7, 10).Value = 90

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 72
